$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The "Longest Repeating Character Replacement" entry (previously just a
# placeholder row with Question/Subject/blind75 filled in) is being fully
# documented like the other finished rows (2-34). That full record is
# inserted at row 35, and the four still-unfinished rows that used to sit at
# 35-38 (Group Anagrams, Top K Fequent Elements, Encode and Decode Strings,
# Valid Palindrome) are pushed down below row 42, becoming rows 43-46.
# ---------------------------------------------------------------------------

# 1) Relocate the four unfinished rows (old 35:38) down to 43:46.
$ws.Range("A35:K38").Copy() | Out-Null
$ws.Range("A43:K43").PasteSpecial(-4104) | Out-Null   # xlPasteAll
$excel.CutCopyMode = 0

# 2) Clear out the old rows 35-39 (35-38 were just copied away above; row 39
#    held the old, incomplete "Longest Repeating Character Replacement" row
#    whose content is being replaced by the fully filled-in version below).
$ws.Range("A35:K39").Clear() | Out-Null

# 3) Fill in the new, complete row 35.
$ws.Range("A35").Value = 424
$ws.Range("B35").Value = "Medium"
$ws.Range("C35").Value = "Longest Repeating Character Replacement"
$ws.Range("D35").Value = "O(n)"
$ws.Range("E35").Value = "Using sliding window technique, check if (cur_len - len_max_repeats) <= k."
$ws.Range("F35").Value = 45510
$ws.Range("G35").Value = "Sliding Window"
$ws.Range("H35").Value = "blind75"
$ws.Range("I35").Value = "Easy"
$ws.Range("J35").Value = "String"
$ws.Range("K35").Value = "https://leetcode.com/problems/longest-repeating-character-replacement/description/"

$ws.Hyperlinks.Add($ws.Range("K35"), "https://leetcode.com/problems/longest-repeating-character-replacement/description/") | Out-Null

# 4) Match row 35's formatting to the other "complete" rows (copy from row 34,
#    which carries the same style pattern) - done after the value/hyperlink
#    writes so it wins over any default styling the hyperlink insertion set.
$ws.Range("A34:K34").Copy() | Out-Null
$ws.Range("A35:K35").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# 5) Widen column K slightly to fit the new, longer hyperlink text.
$ws.Columns("K:K").ColumnWidth = 80

# 6) Reflect the view state: scrolled down a bit further and the whole new
#    row selected (as if freshly typed in and then the row header clicked).
$ws.Range("A35:XFD35").Select() | Out-Null
$win = $excel.ActiveWindow
$win.ScrollRow = 30
$win.ScrollColumn = 1
